$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Update the German translation prompts in column C (values only; the
# cells already carry their styles/positions, so no new cells are added).
$ws.Range("C2").Value = "de_erinnern"
$ws.Range("C3").Value = "de_verlassen"
$ws.Range("C4").Value = "de_beginnen"

# Move the active cell/selection to where the author last left it.
$ws.Range("D10").Select()
